$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 7
$ws.Range("D2").Value = 44301
$ws.Range("I2").Value = "Especial"
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 18000
$ws.Range("P2").Value = 1000

# Row 3 <- original row 31
$ws.Range("D3").Value = 44242
$ws.Range("I3").Value = "Especial"
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 20000
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = 20000
$ws.Range("P3").Value = 1111

# Row 4 <- original row 8
$ws.Range("D4").Value = 44253
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 70
$ws.Range("K4").Value = 18000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 18000
$ws.Range("P4").Value = 1000

# Row 5 <- original row 17
$ws.Range("D5").Value = 44295
$ws.Range("I5").Value = "Especial"
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 16000
$ws.Range("P5").Value = 889

# Row 6 <- original row 5
$ws.Range("D6").Value = 44260
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 70
$ws.Range("K6").Value = 16000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 16000
$ws.Range("P6").Value = 889

# Row 7 <- original row 24
$ws.Range("D7").Value = 44245
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 18000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 18000
$ws.Range("P7").Value = 1000

# Row 8 <- original row 27
$ws.Range("D8").Value = 44236
$ws.Range("I8").Value = "Especial"
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 20000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 20000
$ws.Range("P8").Value = 1111

# Row 9 <- original row 12
$ws.Range("D9").Value = 44320
$ws.Range("I9").Value = "Especial"
$ws.Range("J9").Value = 90
$ws.Range("K9").Value = 17500
$ws.Range("L9").Value = 17500
$ws.Range("M9").Value = 17500
$ws.Range("P9").Value = 972

# Row 10 <- original row 21
$ws.Range("D10").Value = 44326
$ws.Range("I10").Value = "Especial"
$ws.Range("J10").Value = 15
$ws.Range("K10").Value = 18000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 18000
$ws.Range("P10").Value = 1000

# Row 11 <- original row 43
$ws.Range("D11").Value = 44274
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 70
$ws.Range("K11").Value = 16000
$ws.Range("L11").Value = 16000
$ws.Range("M11").Value = 16000
$ws.Range("P11").Value = 889

# Row 12 <- original row 15
$ws.Range("D12").Value = 44309
$ws.Range("I12").Value = "Especial"
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 18000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 18000
$ws.Range("P12").Value = 1000

# Row 13 <- original row 23
$ws.Range("D13").Value = 44252
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 40
$ws.Range("K13").Value = 18000
$ws.Range("L13").Value = 18000
$ws.Range("M13").Value = 18000
$ws.Range("P13").Value = 1000

# Row 15 <- original row 25
$ws.Range("D15").Value = 44243
$ws.Range("I15").Value = "Especial"
$ws.Range("J15").Value = 60
$ws.Range("K15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 20000
$ws.Range("P15").Value = 1111

# Row 16 <- original row 22
$ws.Range("D16").Value = 44292
$ws.Range("I16").Value = "Especial"
$ws.Range("J16").Value = 70
$ws.Range("K16").Value = 17000
$ws.Range("L16").Value = 17000
$ws.Range("M16").Value = 17000
$ws.Range("P16").Value = 944

# Row 17 <- original row 39
$ws.Range("D17").Value = 44299
$ws.Range("I17").Value = "Especial"
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 18000
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = 18000
$ws.Range("P17").Value = 1000

# Row 20 <- original row 40
$ws.Range("D20").Value = 44316
$ws.Range("I20").Value = "Especial"
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = 18000
$ws.Range("L20").Value = 18000
$ws.Range("M20").Value = 18000
$ws.Range("P20").Value = 1000

# Row 21 <- original row 10
$ws.Range("D21").Value = 44364
$ws.Range("I21").Value = "Especial"
$ws.Range("J21").Value = 30
$ws.Range("K21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = 20000
$ws.Range("P21").Value = 1111

# Row 22 <- original row 36
$ws.Range("D22").Value = 44280
$ws.Range("I22").Value = "Especial"
$ws.Range("J22").Value = 40
$ws.Range("K22").Value = 18000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 18000
$ws.Range("P22").Value = 1000

# Row 23 <- original row 16
$ws.Range("D23").Value = 44323
$ws.Range("I23").Value = "Especial"
$ws.Range("J23").Value = 70
$ws.Range("K23").Value = 18000
$ws.Range("L23").Value = 18000
$ws.Range("M23").Value = 18000
$ws.Range("P23").Value = 1000

# Row 24 <- original row 37
$ws.Range("D24").Value = 44267
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 70
$ws.Range("K24").Value = 16000
$ws.Range("L24").Value = 16000
$ws.Range("M24").Value = 16000
$ws.Range("P24").Value = 889

# Row 25 <- original row 28
$ws.Range("D25").Value = 44306
$ws.Range("I25").Value = "Especial"
$ws.Range("J25").Value = 80
$ws.Range("K25").Value = 18000
$ws.Range("L25").Value = 18000
$ws.Range("M25").Value = 18000
$ws.Range("P25").Value = 1000

# Row 26 <- original row 4
$ws.Range("D26").Value = 44271
$ws.Range("I26").Value = "Especial"
$ws.Range("J26").Value = 70
$ws.Range("K26").Value = 18000
$ws.Range("L26").Value = 18000
$ws.Range("M26").Value = 18000
$ws.Range("P26").Value = 1000

# Row 27 <- original row 26
$ws.Range("D27").Value = 44250
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = 18000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 18000
$ws.Range("P27").Value = 1000

# Row 28 <- original row 29
$ws.Range("D28").Value = 44285
$ws.Range("I28").Value = "Especial"
$ws.Range("J28").Value = 70
$ws.Range("K28").Value = 18000
$ws.Range("L28").Value = 18000
$ws.Range("M28").Value = 18000
$ws.Range("P28").Value = 1000

# Row 29 <- original row 34
$ws.Range("D29").Value = 44259
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 70
$ws.Range("K29").Value = 16000
$ws.Range("L29").Value = 16000
$ws.Range("M29").Value = 16000
$ws.Range("P29").Value = 889

# Row 30 <- original row 2
$ws.Range("D30").Value = 44251
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 20
$ws.Range("K30").Value = 18000
$ws.Range("L30").Value = 18000
$ws.Range("M30").Value = 18000
$ws.Range("P30").Value = 1000

# Row 31 <- original row 6
$ws.Range("D31").Value = 44264
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 80
$ws.Range("K31").Value = 16000
$ws.Range("L31").Value = 16000
$ws.Range("M31").Value = 16000
$ws.Range("P31").Value = 889

# Row 32 <- original row 20
$ws.Range("D32").Value = 44232
$ws.Range("I32").Value = "Especial"
$ws.Range("J32").Value = 50
$ws.Range("K32").Value = 22000
$ws.Range("L32").Value = 22000
$ws.Range("M32").Value = 22000
$ws.Range("P32").Value = 1222

# Row 33 <- original row 35
$ws.Range("D33").Value = 44278
$ws.Range("I33").Value = "Especial"
$ws.Range("J33").Value = 70
$ws.Range("K33").Value = 18000
$ws.Range("L33").Value = 18000
$ws.Range("M33").Value = 18000
$ws.Range("P33").Value = 1000

# Row 34 <- original row 42
$ws.Range("D34").Value = 44350
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 20
$ws.Range("K34").Value = 20000
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = 20000
$ws.Range("P34").Value = 1111

# Row 35 <- original row 9
$ws.Range("D35").Value = 44238
$ws.Range("I35").Value = "Especial"
$ws.Range("J35").Value = 50
$ws.Range("K35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("M35").Value = 20000
$ws.Range("P35").Value = 1111

# Row 36 <- original row 3
$ws.Range("D36").Value = 44257
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 60
$ws.Range("K36").Value = 16000
$ws.Range("L36").Value = 16000
$ws.Range("M36").Value = 16000
$ws.Range("P36").Value = 889

# Row 37 <- original row 38
$ws.Range("D37").Value = 44270
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 15
$ws.Range("K37").Value = 16000
$ws.Range("L37").Value = 16000
$ws.Range("M37").Value = 16000
$ws.Range("P37").Value = 889

# Row 38 <- original row 11
$ws.Range("D38").Value = 44302
$ws.Range("I38").Value = "Especial"
$ws.Range("J38").Value = 70
$ws.Range("K38").Value = 18000
$ws.Range("L38").Value = 18000
$ws.Range("M38").Value = 18000
$ws.Range("P38").Value = 1000

# Row 39 <- original row 41
$ws.Range("D39").Value = 44239
$ws.Range("I39").Value = "Especial"
$ws.Range("J39").Value = 60
$ws.Range("K39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("M39").Value = 20000
$ws.Range("P39").Value = 1111

# Row 40 <- original row 33
$ws.Range("D40").Value = 44357
$ws.Range("I40").Value = "Especial"
$ws.Range("J40").Value = 15
$ws.Range("K40").Value = 20000
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = 20000
$ws.Range("P40").Value = 1111

# Row 41 <- original row 30
$ws.Range("D41").Value = 44371
$ws.Range("I41").Value = "Especial"
$ws.Range("J41").Value = 20
$ws.Range("K41").Value = 20000
$ws.Range("L41").Value = 20000
$ws.Range("M41").Value = 20000
$ws.Range("P41").Value = 1111

# Row 42 <- original row 13
$ws.Range("D42").Value = 44365
$ws.Range("I42").Value = "Especial"
$ws.Range("J42").Value = 50
$ws.Range("K42").Value = 20000
$ws.Range("L42").Value = 20000
$ws.Range("M42").Value = 20000
$ws.Range("P42").Value = 1111

# Row 43 <- original row 32
$ws.Range("D43").Value = 44313
$ws.Range("I43").Value = "Especial"
$ws.Range("J43").Value = 80
$ws.Range("K43").Value = 18000
$ws.Range("L43").Value = 18000
$ws.Range("M43").Value = 18000
$ws.Range("P43").Value = 1000
